# This script re-applies refreshed NATMI ligand-receptor (Ntn1-Unc5b) TPM-based
# statistics onto the existing "Ntn1-Unc5b" results sheet, updating the numeric
# columns E:T (ligand/receptor expression stats and derived edge weights) that
# change after recomputing the analysis with the new TPM values. Columns A:D
# (cluster/gene labels) are untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,5).Value = 3.0  # E2 (Ligand-expressing cells)
$ws.Cells.Item(2,6).Value = 1.0  # F2 (Ligand detection rate)
$ws.Cells.Item(2,7).Value = 1.532141  # G2 (Ligand average expression value)
$ws.Cells.Item(2,8).Value = 4.596423  # H2 (Ligand total expression value)
$ws.Cells.Item(2,9).Value = 0.08900664250669833  # I2 (Ligand derived specificity of average expression value)
$ws.Cells.Item(2,10).Value = 0.0890066425066983  # J2 (Ligand derived specificity of total expression value)
$ws.Cells.Item(2,11).Value = 3.0  # K2 (Receptor-expressing cells)
$ws.Cells.Item(2,12).Value = 1.0  # L2 (Receptor detection rate)
$ws.Cells.Item(2,13).Value = 17.33915166666667  # M2 (Receptor average expression value)
$ws.Cells.Item(2,14).Value = 52.01745500000001  # N2 (Receptor total expression value)
$ws.Cells.Item(2,15).Value = 0.7069602198267303  # O2 (Receptor derived specificity of average expression value)
$ws.Cells.Item(2,16).Value = 0.7069602198267302  # P2 (Receptor derived specificity of total expression value)
$ws.Cells.Item(2,17).Value = 26.56602517371834  # Q2 (Edge average expression weight)
$ws.Cells.Item(2,18).Value = 239.094226563465  # R2 (Edge total expression weight)
$ws.Cells.Item(2,19).Value = 0.06292415555257465  # S2 (Edge average expression derived specificity)
$ws.Cells.Item(2,20).Value = 0.06292415555257462  # T2 (Edge total expression derived specificity)

# Row 3
$ws.Cells.Item(3,5).Value = 3.0  # E3 (Ligand-expressing cells)
$ws.Cells.Item(3,6).Value = 1.0  # F3 (Ligand detection rate)
$ws.Cells.Item(3,7).Value = 1.532141  # G3 (Ligand average expression value)
$ws.Cells.Item(3,8).Value = 4.596423  # H3 (Ligand total expression value)
$ws.Cells.Item(3,9).Value = 0.08900664250669833  # I3 (Ligand derived specificity of average expression value)
$ws.Cells.Item(3,10).Value = 0.0890066425066983  # J3 (Ligand derived specificity of total expression value)
$ws.Cells.Item(3,15).Value = 0.2008728953126747  # O3 (Receptor derived specificity of average expression value)
$ws.Cells.Item(3,16).Value = 0.2008728953126747  # P3 (Receptor derived specificity of total expression value)
$ws.Cells.Item(3,17).Value = 7.548365868311667  # Q3 (Edge average expression weight)
$ws.Cells.Item(3,18).Value = 67.93529281480501  # R3 (Edge total expression weight)
$ws.Cells.Item(3,19).Value = 0.01787902198238068  # S3 (Edge average expression derived specificity)
$ws.Cells.Item(3,20).Value = 0.01787902198238067  # T3 (Edge total expression derived specificity)

# Row 4
$ws.Cells.Item(4,5).Value = 3.0  # E4 (Ligand-expressing cells)
$ws.Cells.Item(4,6).Value = 1.0  # F4 (Ligand detection rate)
$ws.Cells.Item(4,7).Value = 1.532141  # G4 (Ligand average expression value)
$ws.Cells.Item(4,8).Value = 4.596423  # H4 (Ligand total expression value)
$ws.Cells.Item(4,9).Value = 0.08900664250669833  # I4 (Ligand derived specificity of average expression value)
$ws.Cells.Item(4,10).Value = 0.0890066425066983  # J4 (Ligand derived specificity of total expression value)
$ws.Cells.Item(4,13).Value = 2.127843333333333  # M4 (Receptor average expression value)
$ws.Cells.Item(4,14).Value = 6.38353  # N4 (Receptor total expression value)
$ws.Cells.Item(4,15).Value = 0.08675745039949621  # O4 (Receptor derived specificity of average expression value)
$ws.Cells.Item(4,16).Value = 0.08675745039949623  # P4 (Receptor derived specificity of total expression value)
$ws.Cells.Item(4,17).Value = 3.260156012576667  # Q4 (Edge average expression weight)
$ws.Cells.Item(4,18).Value = 29.34140411319  # R4 (Edge total expression weight)
$ws.Cells.Item(4,19).Value = 0.007721989372500572  # S4 (Edge average expression derived specificity)
$ws.Cells.Item(4,20).Value = 0.007721989372500571  # T4 (Edge total expression derived specificity)

# Row 5
$ws.Cells.Item(5,5).Value = 3.0  # E5 (Ligand-expressing cells)
$ws.Cells.Item(5,6).Value = 1.0  # F5 (Ligand detection rate)
$ws.Cells.Item(5,7).Value = 1.532141  # G5 (Ligand average expression value)
$ws.Cells.Item(5,8).Value = 4.596423  # H5 (Ligand total expression value)
$ws.Cells.Item(5,9).Value = 0.08900664250669833  # I5 (Ligand derived specificity of average expression value)
$ws.Cells.Item(5,10).Value = 0.0890066425066983  # J5 (Ligand derived specificity of total expression value)
$ws.Cells.Item(5,13).Value = 0.1326736666666667  # M5 (Receptor average expression value)
$ws.Cells.Item(5,14).Value = 0.398021  # N5 (Receptor total expression value)
$ws.Cells.Item(5,15).Value = 0.005409434461098778  # O5 (Receptor derived specificity of average expression value)
$ws.Cells.Item(5,16).Value = 0.005409434461098778  # P5 (Receptor derived specificity of total expression value)
$ws.Cells.Item(5,17).Value = 0.2032747643203333  # Q5 (Edge average expression weight)
$ws.Cells.Item(5,18).Value = 1.829472878883  # R5 (Edge total expression weight)
$ws.Cells.Item(5,19).Value = 0.0004814755992424333  # S5 (Edge average expression derived specificity)
$ws.Cells.Item(5,20).Value = 0.0004814755992424331  # T5 (Edge total expression derived specificity)

# Row 6
$ws.Cells.Item(6,9).Value = 0.6169137955113024  # I6 (Ligand derived specificity of average expression value)
$ws.Cells.Item(6,10).Value = 0.6169137955113023  # J6 (Ligand derived specificity of total expression value)
$ws.Cells.Item(6,11).Value = 3.0  # K6 (Receptor-expressing cells)
$ws.Cells.Item(6,12).Value = 1.0  # L6 (Receptor detection rate)
$ws.Cells.Item(6,13).Value = 17.33915166666667  # M6 (Receptor average expression value)
$ws.Cells.Item(6,14).Value = 52.01745500000001  # N6 (Receptor total expression value)
$ws.Cells.Item(6,15).Value = 0.7069602198267303  # O6 (Receptor derived specificity of average expression value)
$ws.Cells.Item(6,16).Value = 0.7069602198267302  # P6 (Receptor derived specificity of total expression value)
$ws.Cells.Item(6,17).Value = 184.1317339920334  # Q6 (Edge average expression weight)
$ws.Cells.Item(6,18).Value = 1657.1856059283  # R6 (Edge total expression weight)
$ws.Cells.Item(6,19).Value = 0.4361335124888129  # S6 (Edge average expression derived specificity)
$ws.Cells.Item(6,20).Value = 0.4361335124888128  # T6 (Edge total expression derived specificity)

# Row 7
$ws.Cells.Item(7,9).Value = 0.6169137955113024  # I7 (Ligand derived specificity of average expression value)
$ws.Cells.Item(7,10).Value = 0.6169137955113023  # J7 (Ligand derived specificity of total expression value)
$ws.Cells.Item(7,15).Value = 0.2008728953126747  # O7 (Receptor derived specificity of average expression value)
$ws.Cells.Item(7,16).Value = 0.2008728953126747  # P7 (Receptor derived specificity of total expression value)
$ws.Cells.Item(7,19).Value = 0.1239212602626867  # S7 (Edge average expression derived specificity)
$ws.Cells.Item(7,20).Value = 0.1239212602626866  # T7 (Edge total expression derived specificity)

# Row 8
$ws.Cells.Item(8,9).Value = 0.6169137955113024  # I8 (Ligand derived specificity of average expression value)
$ws.Cells.Item(8,10).Value = 0.6169137955113023  # J8 (Ligand derived specificity of total expression value)
$ws.Cells.Item(8,13).Value = 2.127843333333333  # M8 (Receptor average expression value)
$ws.Cells.Item(8,14).Value = 6.38353  # N8 (Receptor total expression value)
$ws.Cells.Item(8,15).Value = 0.08675745039949621  # O8 (Receptor derived specificity of average expression value)
$ws.Cells.Item(8,16).Value = 0.08675745039949623  # P8 (Receptor derived specificity of total expression value)
$ws.Cells.Item(8,17).Value = 22.59646205086667  # Q8 (Edge average expression weight)
$ws.Cells.Item(8,18).Value = 203.3681584578  # R8 (Edge total expression weight)
$ws.Cells.Item(8,19).Value = 0.05352186801483677  # S8 (Edge average expression derived specificity)
$ws.Cells.Item(8,20).Value = 0.05352186801483677  # T8 (Edge total expression derived specificity)

# Row 9
$ws.Cells.Item(9,9).Value = 0.6169137955113024  # I9 (Ligand derived specificity of average expression value)
$ws.Cells.Item(9,10).Value = 0.6169137955113023  # J9 (Ligand derived specificity of total expression value)
$ws.Cells.Item(9,13).Value = 0.1326736666666667  # M9 (Receptor average expression value)
$ws.Cells.Item(9,14).Value = 0.398021  # N9 (Receptor total expression value)
$ws.Cells.Item(9,15).Value = 0.005409434461098778  # O9 (Receptor derived specificity of average expression value)
$ws.Cells.Item(9,16).Value = 0.005409434461098778  # P9 (Receptor derived specificity of total expression value)
$ws.Cells.Item(9,17).Value = 1.408917389273333  # Q9 (Edge average expression weight)
$ws.Cells.Item(9,18).Value = 12.68025650346  # R9 (Edge total expression weight)
$ws.Cells.Item(9,19).Value = 0.003337154744966084  # S9 (Edge average expression derived specificity)
$ws.Cells.Item(9,20).Value = 0.003337154744966083  # T9 (Edge total expression derived specificity)

# Row 10
$ws.Cells.Item(10,7).Value = 4.902263666666666  # G10 (Ligand average expression value)
$ws.Cells.Item(10,8).Value = 14.706791  # H10 (Ligand total expression value)
$ws.Cells.Item(10,9).Value = 0.2847871244569372  # I10 (Ligand derived specificity of average expression value)
$ws.Cells.Item(10,10).Value = 0.2847871244569371  # J10 (Ligand derived specificity of total expression value)
$ws.Cells.Item(10,11).Value = 3.0  # K10 (Receptor-expressing cells)
$ws.Cells.Item(10,12).Value = 1.0  # L10 (Receptor detection rate)
$ws.Cells.Item(10,13).Value = 17.33915166666667  # M10 (Receptor average expression value)
$ws.Cells.Item(10,14).Value = 52.01745500000001  # N10 (Receptor total expression value)
$ws.Cells.Item(10,15).Value = 0.7069602198267303  # O10 (Receptor derived specificity of average expression value)
$ws.Cells.Item(10,16).Value = 0.7069602198267302  # P10 (Receptor derived specificity of total expression value)
$ws.Cells.Item(10,17).Value = 85.00109322632278  # Q10 (Edge average expression weight)
$ws.Cells.Item(10,18).Value = 765.0098390369051  # R10 (Edge total expression weight)
$ws.Cells.Item(10,19).Value = 0.2013331681098987  # S10 (Edge average expression derived specificity)
$ws.Cells.Item(10,20).Value = 0.2013331681098986  # T10 (Edge total expression derived specificity)

# Row 11
$ws.Cells.Item(11,7).Value = 4.902263666666666  # G11 (Ligand average expression value)
$ws.Cells.Item(11,8).Value = 14.706791  # H11 (Ligand total expression value)
$ws.Cells.Item(11,9).Value = 0.2847871244569372  # I11 (Ligand derived specificity of average expression value)
$ws.Cells.Item(11,10).Value = 0.2847871244569371  # J11 (Ligand derived specificity of total expression value)
$ws.Cells.Item(11,15).Value = 0.2008728953126747  # O11 (Receptor derived specificity of average expression value)
$ws.Cells.Item(11,16).Value = 0.2008728953126747  # P11 (Receptor derived specificity of total expression value)
$ws.Cells.Item(11,17).Value = 24.15187619085389  # Q11 (Edge average expression weight)
$ws.Cells.Item(11,18).Value = 217.366885717685  # R11 (Edge total expression weight)
$ws.Cells.Item(11,19).Value = 0.05720601423743599  # S11 (Edge average expression derived specificity)
$ws.Cells.Item(11,20).Value = 0.05720601423743599  # T11 (Edge total expression derived specificity)

# Row 12
$ws.Cells.Item(12,7).Value = 4.902263666666666  # G12 (Ligand average expression value)
$ws.Cells.Item(12,8).Value = 14.706791  # H12 (Ligand total expression value)
$ws.Cells.Item(12,9).Value = 0.2847871244569372  # I12 (Ligand derived specificity of average expression value)
$ws.Cells.Item(12,10).Value = 0.2847871244569371  # J12 (Ligand derived specificity of total expression value)
$ws.Cells.Item(12,13).Value = 2.127843333333333  # M12 (Receptor average expression value)
$ws.Cells.Item(12,14).Value = 6.38353  # N12 (Receptor total expression value)
$ws.Cells.Item(12,15).Value = 0.08675745039949621  # O12 (Receptor derived specificity of average expression value)
$ws.Cells.Item(12,16).Value = 0.08675745039949623  # P12 (Receptor derived specificity of total expression value)
$ws.Cells.Item(12,17).Value = 10.43124906135889  # Q12 (Edge average expression weight)
$ws.Cells.Item(12,18).Value = 93.88124155223  # R12 (Edge total expression weight)
$ws.Cells.Item(12,19).Value = 0.02470740482448788  # S12 (Edge average expression derived specificity)
$ws.Cells.Item(12,20).Value = 0.02470740482448788  # T12 (Edge total expression derived specificity)

# Row 13
$ws.Cells.Item(13,7).Value = 4.902263666666666  # G13 (Ligand average expression value)
$ws.Cells.Item(13,8).Value = 14.706791  # H13 (Ligand total expression value)
$ws.Cells.Item(13,9).Value = 0.2847871244569372  # I13 (Ligand derived specificity of average expression value)
$ws.Cells.Item(13,10).Value = 0.2847871244569371  # J13 (Ligand derived specificity of total expression value)
$ws.Cells.Item(13,13).Value = 0.1326736666666667  # M13 (Receptor average expression value)
$ws.Cells.Item(13,14).Value = 0.398021  # N13 (Receptor total expression value)
$ws.Cells.Item(13,15).Value = 0.005409434461098778  # O13 (Receptor derived specificity of average expression value)
$ws.Cells.Item(13,16).Value = 0.005409434461098778  # P13 (Receptor derived specificity of total expression value)
$ws.Cells.Item(13,17).Value = 0.6504012956234444  # Q13 (Edge average expression weight)
$ws.Cells.Item(13,18).Value = 5.853611660610999  # R13 (Edge total expression weight)
$ws.Cells.Item(13,19).Value = 0.001540537285114583  # S13 (Edge average expression derived specificity)
$ws.Cells.Item(13,20).Value = 0.001540537285114582  # T13 (Edge total expression derived specificity)

# Row 14
$ws.Cells.Item(14,5).Value = 2.0  # E14 (Ligand-expressing cells)
$ws.Cells.Item(14,6).Value = 0.6666666666666666  # F14 (Ligand detection rate)
$ws.Cells.Item(14,7).Value = 0.159958  # G14 (Ligand average expression value)
$ws.Cells.Item(14,8).Value = 0.479874  # H14 (Ligand total expression value)
$ws.Cells.Item(14,9).Value = 0.009292437525062282  # I14 (Ligand derived specificity of average expression value)
$ws.Cells.Item(14,10).Value = 0.00929243752506228  # J14 (Ligand derived specificity of total expression value)
$ws.Cells.Item(14,11).Value = 3.0  # K14 (Receptor-expressing cells)
$ws.Cells.Item(14,12).Value = 1.0  # L14 (Receptor detection rate)
$ws.Cells.Item(14,13).Value = 17.33915166666667  # M14 (Receptor average expression value)
$ws.Cells.Item(14,14).Value = 52.01745500000001  # N14 (Receptor total expression value)
$ws.Cells.Item(14,15).Value = 0.7069602198267303  # O14 (Receptor derived specificity of average expression value)
$ws.Cells.Item(14,16).Value = 0.7069602198267302  # P14 (Receptor derived specificity of total expression value)
$ws.Cells.Item(14,17).Value = 2.773536022296667  # Q14 (Edge average expression weight)
$ws.Cells.Item(14,18).Value = 24.96182420067  # R14 (Edge total expression weight)
$ws.Cells.Item(14,19).Value = 0.006569383675444189  # S14 (Edge average expression derived specificity)
$ws.Cells.Item(14,20).Value = 0.006569383675444187  # T14 (Edge total expression derived specificity)

# Row 15
$ws.Cells.Item(15,5).Value = 2.0  # E15 (Ligand-expressing cells)
$ws.Cells.Item(15,6).Value = 0.6666666666666666  # F15 (Ligand detection rate)
$ws.Cells.Item(15,7).Value = 0.159958  # G15 (Ligand average expression value)
$ws.Cells.Item(15,8).Value = 0.479874  # H15 (Ligand total expression value)
$ws.Cells.Item(15,9).Value = 0.009292437525062282  # I15 (Ligand derived specificity of average expression value)
$ws.Cells.Item(15,10).Value = 0.00929243752506228  # J15 (Ligand derived specificity of total expression value)
$ws.Cells.Item(15,15).Value = 0.2008728953126747  # O15 (Receptor derived specificity of average expression value)
$ws.Cells.Item(15,16).Value = 0.2008728953126747  # P15 (Receptor derived specificity of total expression value)
$ws.Cells.Item(15,17).Value = 0.7880616128433334  # Q15 (Edge average expression weight)
$ws.Cells.Item(15,18).Value = 7.092554515590001  # R15 (Edge total expression weight)
$ws.Cells.Item(15,19).Value = 0.001866598830171406  # S15 (Edge average expression derived specificity)
$ws.Cells.Item(15,20).Value = 0.001866598830171405  # T15 (Edge total expression derived specificity)

# Row 16
$ws.Cells.Item(16,5).Value = 2.0  # E16 (Ligand-expressing cells)
$ws.Cells.Item(16,6).Value = 0.6666666666666666  # F16 (Ligand detection rate)
$ws.Cells.Item(16,7).Value = 0.159958  # G16 (Ligand average expression value)
$ws.Cells.Item(16,8).Value = 0.479874  # H16 (Ligand total expression value)
$ws.Cells.Item(16,9).Value = 0.009292437525062282  # I16 (Ligand derived specificity of average expression value)
$ws.Cells.Item(16,10).Value = 0.00929243752506228  # J16 (Ligand derived specificity of total expression value)
$ws.Cells.Item(16,13).Value = 2.127843333333333  # M16 (Receptor average expression value)
$ws.Cells.Item(16,14).Value = 6.38353  # N16 (Receptor total expression value)
$ws.Cells.Item(16,15).Value = 0.08675745039949621  # O16 (Receptor derived specificity of average expression value)
$ws.Cells.Item(16,16).Value = 0.08675745039949623  # P16 (Receptor derived specificity of total expression value)
$ws.Cells.Item(16,17).Value = 0.3403655639133333  # Q16 (Edge average expression weight)
$ws.Cells.Item(16,18).Value = 3.06329007522  # R16 (Edge total expression weight)
$ws.Cells.Item(16,19).Value = 0.0008061881876710083  # S16 (Edge average expression derived specificity)
$ws.Cells.Item(16,20).Value = 0.0008061881876710083  # T16 (Edge total expression derived specificity)

# Row 17
$ws.Cells.Item(17,5).Value = 2.0  # E17 (Ligand-expressing cells)
$ws.Cells.Item(17,6).Value = 0.6666666666666666  # F17 (Ligand detection rate)
$ws.Cells.Item(17,7).Value = 0.159958  # G17 (Ligand average expression value)
$ws.Cells.Item(17,8).Value = 0.479874  # H17 (Ligand total expression value)
$ws.Cells.Item(17,9).Value = 0.009292437525062282  # I17 (Ligand derived specificity of average expression value)
$ws.Cells.Item(17,10).Value = 0.00929243752506228  # J17 (Ligand derived specificity of total expression value)
$ws.Cells.Item(17,13).Value = 0.1326736666666667  # M17 (Receptor average expression value)
$ws.Cells.Item(17,14).Value = 0.398021  # N17 (Receptor total expression value)
$ws.Cells.Item(17,15).Value = 0.005409434461098778  # O17 (Receptor derived specificity of average expression value)
$ws.Cells.Item(17,16).Value = 0.005409434461098778  # P17 (Receptor derived specificity of total expression value)
$ws.Cells.Item(17,17).Value = 0.02122221437266666  # Q17 (Edge average expression weight)
$ws.Cells.Item(17,18).Value = 0.190999929354  # R17 (Edge total expression weight)
$ws.Cells.Item(17,19).Value = 0.00005026683177567936  # S17 (Edge average expression derived specificity)
$ws.Cells.Item(17,20).Value = 0.00005026683177567933  # T17 (Edge total expression derived specificity)
